# Fruta / hortaliza, semanal
#
# Weekly refresh of the "Feria Lagunitas de Puerto Montt - Espinaca" data:
# a new weekly observation (fecha serial 45223) is inserted as row 90,
# pushing the previously existing rows 90-95 down to rows 91-96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 90 (shifts rows 90:95 down to 91:96,
# inheriting the formatting - including the date style - of the row above).
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the latest weekly record.
$ws.Cells.Item(90, 1).Value  = 4
$ws.Cells.Item(90, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(90, 3).Value  = "Los Lagos"
$ws.Cells.Item(90, 4).Value  = 45223
$ws.Cells.Item(90, 5).Value  = 10
$ws.Cells.Item(90, 6).Value  = 100112012
$ws.Cells.Item(90, 7).Value  = "Espinaca"
$ws.Cells.Item(90, 8).Value  = "Sin especificar"
$ws.Cells.Item(90, 9).Value  = "Primera"
$ws.Cells.Item(90, 10).Value = 35
$ws.Cells.Item(90, 11).Value = 13000
$ws.Cells.Item(90, 12).Value = 13000
$ws.Cells.Item(90, 13).Value = 13000
$ws.Cells.Item(90, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(90, 15).Value = "Región Metropolitana"
$ws.Cells.Item(90, 16).Value = 1300
$ws.Cells.Item(90, 17).Value = 10
$ws.Cells.Item(90, 18).Value = "Hortaliza"
